$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style/format from G1 (bold, bordered, centered) to H1, then set its value.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill the new "Save" column values for rows 2-9.
$saveValues = @(0, 0, 1, 1, 1, 0, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
